$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell writes are ordered to reproduce the author's original shared-string
# insertion sequence (new unique strings get appended to sharedStrings.xml
# in the order they are first written).
$ws.Range("F7").Value  = 'Triangle ="Red",  3,  3,  5'
$ws.Range("E13").Value = 'Triangle ="Red",  5, 5 ,  6'
$ws.Range("F8").Value  = 'Triangle ="   ",  5,  5,  6'
$ws.Range("F9").Value  = 'Triangle ="Red", five,  5,  6'
$ws.Range("F10").Value = 'Triangle ="Red", 5,  five,  6'
$ws.Range("F11").Value = 'Triangle ="Red",  5,  5,  six'
$ws.Range("E12").Value = 'Triangle ="Red",  5,  5,  6'
$ws.Range("E14").Value = 'Triangle ="Red",  5,  5 ,  6'
$ws.Range("G13").Value = "Area  12"
$ws.Range("G14").Value = "Perimeter  16"

# These reuse already-existing shared strings, so ordering relative to the
# block above does not affect the resulting sharedStrings.xml layout.
$ws.Range("F12").Value = "None"
$ws.Range("F13").Value = "None"
$ws.Range("F14").Value = "None"
$ws.Range("G12").Value = "The shape color is Red. This triangle has three sides with lengths of 5 ,  5 and 6 centimeters."

$ws.Range("E14").Select()
